$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 608
$ws.Range("I2").Value = 219.33333
$ws.Range("K2").Value = 219.33333
$ws.Range("M2").Value = -106.33333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 7059.476
$ws.Range("I76").Value = 5270.4165
$ws.Range("K76").Value = 5270.4165
$ws.Range("M76").Value = -4955.4165

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 7059.476
$ws.Range("I79").Value = 5270.4165
$ws.Range("K79").Value = 5270.4165
$ws.Range("M79").Value = -4178.4165

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 2490.4
$ws.Range("I94").Value = 2988
$ws.Range("K94").Value = 2988
$ws.Range("M94").Value = -2537

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3728.5334
$ws.Range("I100").Value = 2083.25
$ws.Range("K100").Value = 2083.25
$ws.Range("M100").Value = -1542.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1975.4286
$ws.Range("J112").Value = 2010.2
$ws.Range("L112").Value = 6030.6
$ws.Range("N112").Value = -8246.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 650.3570999999999
$ws.Range("I132").Value = 548.4314000000001
$ws.Range("J132").Value = 1690
$ws.Range("K132").Value = 1645.2942
$ws.Range("L132").Value = 5070
$ws.Range("M132").Value = 884.7057999999997
$ws.Range("N132").Value = -10130

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 773.82355
$ws.Range("I135").Value = 697.25
$ws.Range("J135").Value = 1999
$ws.Range("K135").Value = 6275.25
$ws.Range("L135").Value = 17991
$ws.Range("M135").Value = -3740.25
$ws.Range("N135").Value = -23061

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 13891548
$ws.Range("I137").Value = 52632932
$ws.Range("K137").Value = 157898796
$ws.Range("M137").Value = -157896246

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2886.9785
$ws.Range("J138").Value = 3326.0747
$ws.Range("L138").Value = 9978.224099999999
$ws.Range("N138").Value = -20258.2241

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3840.7847
$ws.Range("I32").Value = 2535.377
$ws.Range("K32").Value = 2535.377
$ws.Range("M32").Value = -2248.377

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 22226466
$ws.Range("I74").Value = 25644204
$ws.Range("J74").Value = 11164
$ws.Range("K74").Value = 25644204
$ws.Range("L74").Value = 11164
$ws.Range("M74").Value = -25643330
$ws.Range("N74").Value = -12912

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 22226466
$ws.Range("I77").Value = 25644204
$ws.Range("J77").Value = 11164
$ws.Range("K77").Value = 128221020
$ws.Range("L77").Value = 55820
$ws.Range("M77").Value = -128216652
$ws.Range("N77").Value = -64556

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1331.091
$ws.Range("I97").Value = 1330.9
$ws.Range("J97").Value = 1333
$ws.Range("K97").Value = 1330.9
$ws.Range("L97").Value = 1333
$ws.Range("M97").Value = -834.9000000000001
$ws.Range("N97").Value = -2325

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 7401.0527
$ws.Range("I110").Value = 5247.1333
$ws.Range("J110").Value = 15478.25
$ws.Range("K110").Value = 5247.1333
$ws.Range("L110").Value = 15478.25
$ws.Range("M110").Value = -3202.1333
$ws.Range("N110").Value = -19568.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 6249.9165
$ws.Range("I122").Value = 6997.8
$ws.Range("K122").Value = 20993.4
$ws.Range("M122").Value = -18543.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 11666.667
$ws.Range("J124").Value = 11666.667
$ws.Range("L124").Value = 11666.667
$ws.Range("N124").Value = -21486.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 600238.3
$ws.Range("J125").Value = 600238.3
$ws.Range("L125").Value = 600238.3
$ws.Range("N125").Value = -610078.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H131").Value = 77699.5
$ws.Range("J131").Value = 77699.5
$ws.Range("L131").Value = 77699.5
$ws.Range("N131").Value = -87779.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 414
$ws.Range("I22").Value = 218.66667
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 218.66667
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -45.66667000000001
$ws.Range("N22").Value = -1346

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 6653.8213
$ws.Range("J105").Value = 14101.363
$ws.Range("L105").Value = 14101.363
$ws.Range("N105").Value = -17595.363

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 62614.08
$ws.Range("J132").Value = 62614.08
$ws.Range("L132").Value = 62614.08
$ws.Range("N132").Value = -72734.08

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 6133.6665
$ws.Range("I122").Value = 2423
$ws.Range("J122").Value = 11699.667
$ws.Range("K122").Value = 7269
$ws.Range("L122").Value = 35099.001
$ws.Range("M122").Value = -4819
$ws.Range("N122").Value = -39999.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 4833
$ws.Range("I14").Value = 4833
$ws.Range("K14").Value = 14499
$ws.Range("M14").Value = -14326

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 3286.5715
$ws.Range("I86").Value = 1200
$ws.Range("J86").Value = 4121.2
$ws.Range("K86").Value = 3600
$ws.Range("L86").Value = 12363.6
$ws.Range("M86").Value = -2414
$ws.Range("N86").Value = -14735.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 3286.5715
$ws.Range("I89").Value = 1200
$ws.Range("J89").Value = 4121.2
$ws.Range("K89").Value = 10800
$ws.Range("L89").Value = 37090.8
$ws.Range("M89").Value = -4872
$ws.Range("N89").Value = -48946.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4527.077
$ws.Range("I102").Value = 1207.5
$ws.Range("J102").Value = 6002.4443
$ws.Range("K102").Value = 1207.5
$ws.Range("L102").Value = 6002.4443
$ws.Range("M102").Value = 414.5
$ws.Range("N102").Value = -9246.444299999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 13340.131
$ws.Range("J122").Value = 10217.2
$ws.Range("L122").Value = 30651.6
$ws.Range("N122").Value = -35551.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7603.5
$ws.Range("I40").Value = 6725.353
$ws.Range("K40").Value = 6725.353
$ws.Range("M40").Value = -6589.353

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 9710.143
$ws.Range("J61").Value = 11500.5
$ws.Range("L61").Value = 11500.5
$ws.Range("N61").Value = -11904.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 17689.062
$ws.Range("I100").Value = 17285.285
$ws.Range("K100").Value = 17285.285
$ws.Range("M100").Value = -16744.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 9710.143
$ws.Range("J113").Value = 11500.5
$ws.Range("L113").Value = 11500.5
$ws.Range("N113").Value = -15840.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6188.875
$ws.Range("I126").Value = 6859.4287
$ws.Range("J126").Value = 1495
$ws.Range("K126").Value = 20578.2861
$ws.Range("L126").Value = 4485
$ws.Range("M126").Value = -18108.2861
$ws.Range("N126").Value = -9425

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2084.85
$ws.Range("I132").Value = 1650.5807
$ws.Range("K132").Value = 4951.742099999999
$ws.Range("M132").Value = -2421.742099999999
